$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: update title (D9) and link (E9)
$ws.Range("D9").Value = "[대학원] 폴란드 DS 대학원 갔던 학생에게 받은 교육 후기 + DS질문 + 답변"
$ws.Range("E9").Value = "https://pdsi.pabii.com/poland-ds-grad-school-review/#utm_source=rss&utm_medium=rss&utm_campaign=poland-ds-grad-school-review"

# Row 28: update title (D28) and link (E28)
$ws.Range("D28").Value = "ROS2 ::  Moveit2 Python API에 대해 알아보기 (Rolling)"
$ws.Range("E28").Value = "https://ropiens.tistory.com/220"

# Row 51: update title (D51) and link (E51)
$ws.Range("D51").Value = "[plotly.js] 차트, 플롯의 배경색 설정하는 방법"
$ws.Range("E51").Value = "https://bskyvision.com/entry/plotlyjs-%EC%B0%A8%ED%8A%B8-%ED%94%8C%EB%A1%AF%EC%9D%98-%EB%B0%B0%EA%B2%BD%EC%83%89-%EC%84%A4%EC%A0%95%ED%95%98%EB%8A%94-%EB%B0%A9%EB%B2%95"
